$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.690303
$ws.Range("H2").Value = 17.070909
$ws.Range("I2").Value = 0.3759010823723209
$ws.Range("J2").Value = 0.3759010823723208
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 10.91915846528567
$ws.Range("R2").Value = 98.272426187571
$ws.Range("S2").Value = 0.002452385864604829
$ws.Range("T2").Value = 0.002452385864604829
$ws.Range("G3").Value = 5.690303
$ws.Range("H3").Value = 17.070909
$ws.Range("I3").Value = 0.3759010823723209
$ws.Range("J3").Value = 0.3759010823723208
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 1031.585877069241
$ws.Range("R3").Value = 9284.272893623167
$ws.Range("S3").Value = 0.2316887909533966
$ws.Range("T3").Value = 0.2316887909533965
$ws.Range("G4").Value = 5.690303
$ws.Range("H4").Value = 17.070909
$ws.Range("I4").Value = 0.3759010823723209
$ws.Range("J4").Value = 0.3759010823723208
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 165.2900570213837
$ws.Range("R4").Value = 1487.610513192453
$ws.Range("S4").Value = 0.03712328204482768
$ws.Range("T4").Value = 0.03712328204482768
$ws.Range("G5").Value = 5.690303
$ws.Range("H5").Value = 17.070909
$ws.Range("I5").Value = 0.3759010823723209
$ws.Range("J5").Value = 0.3759010823723208
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 465.8907433217827
$ws.Range("R5").Value = 4193.016689896044
$ws.Range("S5").Value = 0.1046366235094918
$ws.Range("T5").Value = 0.1046366235094918
$ws.Range("I6").Value = 0.2572281411562661
$ws.Range("J6").Value = 0.2572281411562661
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 7.471951975477888
$ws.Range("R6").Value = 67.24756777930099
$ws.Range("S6").Value = 0.001678161322040005
$ws.Range("T6").Value = 0.001678161322040005
$ws.Range("I7").Value = 0.2572281411562661
$ws.Range("J7").Value = 0.2572281411562661
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.1585440420856667
$ws.Range("T7").Value = 0.1585440420856667
$ws.Range("I8").Value = 0.2572281411562661
$ws.Range("J8").Value = 0.2572281411562661
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 113.1075597093159
$ws.Range("R8").Value = 1017.968037383843
$ws.Range("S8").Value = 0.02540336615618629
$ws.Range("T8").Value = 0.0254033661561863
$ws.Range("I9").Value = 0.2572281411562661
$ws.Range("J9").Value = 0.2572281411562661
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 318.8078340457515
$ws.Range("R9").Value = 2869.270506411764
$ws.Range("S9").Value = 0.07160257159237306
$ws.Range("T9").Value = 0.07160257159237306
$ws.Range("G10").Value = 4.648693333333333
$ws.Range("H10").Value = 13.94608
$ws.Range("I10").Value = 0.3070924088957991
$ws.Range("J10").Value = 0.307092408895799
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 8.920407079057776
$ws.Range("R10").Value = 80.28366371151998
$ws.Range("S10").Value = 0.002003476760297189
$ws.Range("T10").Value = 0.002003476760297188
$ws.Range("G11").Value = 4.648693333333333
$ws.Range("H11").Value = 13.94608
$ws.Range("I11").Value = 0.3070924088957991
$ws.Range("J11").Value = 0.307092408895799
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 842.7541362019912
$ws.Range("R11").Value = 7584.78722581792
$ws.Range("S11").Value = 0.1892781698818349
$ws.Range("T11").Value = 0.1892781698818349
$ws.Range("G12").Value = 4.648693333333333
$ws.Range("H12").Value = 13.94608
$ws.Range("I12").Value = 0.3070924088957991
$ws.Range("J12").Value = 0.307092408895799
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 135.0337207248178
$ws.Range("R12").Value = 1215.30348652336
$ws.Range("S12").Value = 0.03032786720728992
$ws.Range("T12").Value = 0.03032786720728992
$ws.Range("G13").Value = 4.648693333333333
$ws.Range("H13").Value = 13.94608
$ws.Range("I13").Value = 0.3070924088957991
$ws.Range("J13").Value = 0.307092408895799
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 380.6094671130311
$ws.Range("R13").Value = 3425.48520401728
$ws.Range("S13").Value = 0.08548289504637704
$ws.Range("T13").Value = 0.08548289504637703
$ws.Range("G14").Value = 0.9049109999999999
$ws.Range("H14").Value = 2.714733
$ws.Range("I14").Value = 0.05977836757561403
$ws.Range("J14").Value = 0.05977836757561403
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 1.736439449003
$ws.Range("R14").Value = 15.627955041027
$ws.Range("S14").Value = 0.0003899952155668021
$ws.Range("T14").Value = 0.0003899952155668021
$ws.Range("G15").Value = 0.9049109999999999
$ws.Range("H15").Value = 2.714733
$ws.Range("I15").Value = 0.05977836757561403
$ws.Range("J15").Value = 0.05977836757561403
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 164.049859489838
$ws.Range("R15").Value = 1476.448735408542
$ws.Range("S15").Value = 0.03684474016769037
$ws.Range("T15").Value = 0.03684474016769037
$ws.Range("G16").Value = 0.9049109999999999
$ws.Range("H16").Value = 2.714733
$ws.Range("I16").Value = 0.05977836757561403
$ws.Range("J16").Value = 0.05977836757561403
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 26.285558218829
$ws.Range("R16").Value = 236.570023969461
$ws.Range("S16").Value = 0.005903598855538458
$ws.Range("T16").Value = 0.005903598855538459
$ws.Range("G17").Value = 0.9049109999999999
$ws.Range("H17").Value = 2.714733
$ws.Range("I17").Value = 0.05977836757561403
$ws.Range("J17").Value = 0.05977836757561403
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 74.08914049569199
$ws.Range("R17").Value = 666.802264461228
$ws.Range("S17").Value = 0.01664003333681839
$ws.Range("T17").Value = 0.01664003333681839
